$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# YDS sheet: append newly logged week's yardage samples to the running,
# space-separated play-log strings.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = $ws.Range("B2").Value2 + " 4 2 4 4 3 5 -1 5 2 2 1 3 2 6 7 2 2 5 6 14 4 5 7 1 0 -1 7 7 13 1 -2 10 1 2"
$ws.Range("B3").Value = $ws.Range("B3").Value2 + " 5 8 7 28 -1 10 3 11 -1 8 9 11 8 25 2 24 24 5"
$ws.Range("C2").Value = $ws.Range("C2").Value2 + " 6 0 2 4 5 5 6 5 1 -1 5 8 4 3 5 22 2 12 1 8 2 4 1 1 3"
$ws.Range("C3").Value = $ws.Range("C3").Value2 + " 6 2 6 17 2 2 9 8 4 19 3 56 4 4 15"

# ---------------------------------------------------------------------------
# OFF sheet: season-to-date offensive totals updated for the new week.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("OFF")

$ws.Range("C2").Value = 175
$ws.Range("F2").Value = 58
$ws.Range("G2").Value = 54
$ws.Range("J2").Value = 25
$ws.Range("L2").Value = 256
$ws.Range("M2").Value = 165
$ws.Range("O2").Value = 23
$ws.Range("Q2").Value = 505

$ws.Range("B3").Value = 9
$ws.Range("C3").Value = 147
$ws.Range("E3").Value = 30
$ws.Range("F3").Value = 82
$ws.Range("G3").Value = 26
$ws.Range("H3").Value = 33
$ws.Range("I3").Value = 60
$ws.Range("J3").Value = 38
$ws.Range("N3").Value = 22

# ---------------------------------------------------------------------------
# DEF sheet: season-to-date defensive totals updated for the new week.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("DEF")

$ws.Range("C2").Value = 157
$ws.Range("F2").Value = 55
$ws.Range("G2").Value = 43
$ws.Range("H2").Value = 4
$ws.Range("I2").Value = 7
$ws.Range("J2").Value = 29
$ws.Range("L2").Value = 266
$ws.Range("M2").Value = 166
$ws.Range("Q2").Value = 484

$ws.Range("C3").Value = 132
$ws.Range("E3").Value = 29
$ws.Range("F3").Value = 70
$ws.Range("G3").Value = 36
$ws.Range("H3").Value = 26
$ws.Range("I3").Value = 46
$ws.Range("J3").Value = 53
$ws.Range("N3").Value = 24

# ---------------------------------------------------------------------------
# ST sheet: special-teams totals + play-log strings updated for the new week.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ST")

$ws.Range("B2").Value = 65
$ws.Range("D2").Value = 54
$ws.Range("F2").Value = 254
$ws.Range("G2").Value = 247
$ws.Range("N2").Value = 56
$ws.Range("O2").Value = 30

$ws.Range("B3").Value = 54
$ws.Range("D3").Value = $ws.Range("D3").Value2 + " 51 57 44 48 63 46"

$ws.Range("B4").Value = $ws.Range("B4").Value2 + " 71"
$ws.Range("D4").Value = $ws.Range("D4").Value2 + " 12 0 0 5 0 5"

$ws.Range("B5").Value = $ws.Range("B5").Value2 + " 23"
$ws.Range("D5").Value = $ws.Range("D5").Value2 + " 0 0 0 0 0 12 18"

# ---------------------------------------------------------------------------
# TURNS sheet: turnover counts updated for the new week.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("TURNS")

$ws.Range("D2").Value = 7
$ws.Range("E3").Value = 5

# ---------------------------------------------------------------------------
# PEN sheet: penalty counts updated for the new week.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("PEN")

$ws.Range("B3").Value = 14
$ws.Range("D3").Value = 3
